# Apply "started integrating cometchat for user to user messaging" edit:
# Populate the "Marketing Event" column (G) with additional entries for
# several months, including two brand-new marketing events.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New marketing events not previously present in the sheet.
$ws.Range("G4").Value = "post on product hunt "
$ws.Range("G27").Value = "Spread in Entrepreneur Magazine"

# Reuse of existing marketing event text in additional months.
$ws.Range("G14").Value = "bump up # facebbok ads"
$ws.Range("G18").Value = "articles on dev sites about DA"
$ws.Range("G20").Value = "articles on dev sites about DA"
$ws.Range("G22").Value = "articles on dev sites about DA"
$ws.Range("G24").Value = "articles on dev sites about DA"
$ws.Range("G33").Value = "targeted tv ads"

# G25 previously held "targeted tv ads"; swap it to "bump up # facebbok ads".
$ws.Range("G25").Value = "bump up # facebbok ads"

# Move the active selection to I16, matching the author's final cursor spot.
$ws.Range("I16").Select()
